# Generate Report for Handoff
# - Overview sheet: "Latest HO Xliff Generate Date" (col G, rows 4-7) timestamp advances
#   from 18:32:34 to 18:32:50
# - zh-cn sheet: Priority (col E, rows 4-7) moves from "low" to "ht"; the
#   "Latest Handoff Datetime" (col H, rows 4-7) timestamp advances from
#   18:32:29 to 18:32:44
# - de-de sheet: Priority (col E, rows 4-7) moves from "low" to "ht"; the
#   "Latest Handoff Datetime" (col H, rows 4-7) timestamp advances from
#   18:32:34 to 18:32:50 (shares the same value as the Overview sheet)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-18 18:32:50"

    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-18 18:32:44"

    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-18 18:32:50"
}
